# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder country names (sharedStrings effect) ---
# Swap Malta (row 103) and Nigeria (row 104): Nigeria now comes first.
$ws.Range("A103").Value = "Nigeria"
$ws.Range("A104").Value = "Malta"

# Move Santa Lucia ahead of Mongolia/Namibia:
# old order rows 170-172: Mongolia, Namibia, Santa Lucia
# new order rows 170-172: Santa Lucia, Mongolia, Namibia
$ws.Range("A170").Value = "Santa Lucia"
$ws.Range("A171").Value = "Mongolia"
$ws.Range("A172").Value = "Namibia"

# --- Update "last updated" footer text ---
$ws.Range("A1").Value = "Datos actualizados a 4 de Abril de 2020 a las 23:52"

# --- Update numeric stats ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 307689
$ws.Range("C4").Value = 30528
$ws.Range("E4").Value = 284626
$ws.Range("G4").Value = 973
$ws.Range("H4").Value = 8377

# Row 5: España
$ws.Range("B5").Value = 126168
$ws.Range("C5").Value = 6969
$ws.Range("E5").Value = 80002
$ws.Range("G5").Value = 749
$ws.Range("H5").Value = 11947

# Row 7: Alemania
$ws.Range("B7").Value = 95648
$ws.Range("C7").Value = 4489
$ws.Range("E7").Value = 67820
$ws.Range("G7").Value = 153
$ws.Range("H7").Value = 1428

# Row 87: Uruguay
$ws.Range("E87").Value = 295
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 5

# Row 90: Albania
$ws.Range("E90").Value = 214
$ws.Range("G90").Value = 3
$ws.Range("H90").Value = 20

# Row 103 now holds Nigeria's stats
$ws.Range("B103").Value = 214
$ws.Range("C103").Value = 4
$ws.Range("D103").Value = 25
$ws.Range("E103").Value = 185
$ws.Range("F103").Value = 2
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 4

# Row 104 now holds Malta's stats
$ws.Range("B104").Value = 213
$ws.Range("C104").Value = 11
$ws.Range("D104").Value = 2
$ws.Range("E104").Value = 211
$ws.Range("F104").Value = 3
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 0

# Row 170 now holds Santa Lucia's stats
$ws.Range("B170").Value = 14
$ws.Range("C170").Value = 1
$ws.Range("D170").Value = 1
$ws.Range("E170").Value = 13
$ws.Range("F170").Value = 0
$ws.Range("G170").Value = 0
$ws.Range("H170").Value = 0

# Row 171 now holds Mongolia's stats
$ws.Range("B171").Value = 14
$ws.Range("C171").Value = 0
$ws.Range("D171").Value = 2
$ws.Range("E171").Value = 12
$ws.Range("F171").Value = 0
$ws.Range("G171").Value = 0
$ws.Range("H171").Value = 0

# Row 172 now holds Namibia's stats
$ws.Range("B172").Value = 14
$ws.Range("C172").Value = 0
$ws.Range("D172").Value = 3
$ws.Range("E172").Value = 11
$ws.Range("F172").Value = 0
$ws.Range("G172").Value = 0
$ws.Range("H172").Value = 0
